$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "18-012180-417"
$ws.Range("A2").Value = "18-006979-297"

$ws.Range("A1:A2").Style = "Normal"
$ws.Range("A1:A2").Interior.Color = 65535

$ws.Range("A2").Select()
